# Update "想去人数" (desired-attendee count) figures across all four sheets
# to reflect newer scrape numbers (column F), per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 147
    3  = 197
    6  = 1301
    7  = 66
    10 = 444
    11 = 814
    12 = 205
    13 = 746
    14 = 310
    15 = 462
    16 = 91
    17 = 1049
    18 = 492
    20 = 407
    21 = 101
    22 = 218
    23 = 29
    24 = 54
    26 = 433
    27 = 281
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    7  = 291
    8  = 88
    12 = 142
    14 = 8
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 354

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 354
    4  = 147
    5  = 197
    8  = 1301
    10 = 66
    16 = 291
    17 = 444
    18 = 814
    19 = 205
    20 = 746
    21 = 310
    22 = 462
    23 = 91
    24 = 1049
    25 = 492
    26 = 88
    29 = 407
    31 = 101
    33 = 218
    34 = 29
    35 = 54
    36 = 142
    39 = 8
    41 = 433
    42 = 281
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
